$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "56.917.73"
Set-TextValue $ws.Range("E2") "  +0.12%  "
Set-TextValue $ws.Range("D3") "3.027.53"
Set-TextValue $ws.Range("E3") "  +1.83%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "511.74"
Set-TextValue $ws.Range("E5") "  +2.87%  "
Set-TextValue $ws.Range("D6") "140.60"
Set-TextValue $ws.Range("E6") "  +2.59%  "
Set-TextValue $ws.Range("E7") "  -0.03%  "
Set-TextValue $ws.Range("D8") "0.438"
Set-TextValue $ws.Range("E8") "  +2.58%  "
Set-TextValue $ws.Range("E9") "  -2.23%  "
Set-TextValue $ws.Range("E10") "  +1.22%  "
Set-TextValue $ws.Range("D11") "0.375"
Set-TextValue $ws.Range("E11") "  +5.35%  "
Set-TextValue $ws.Range("D12") "3.560.08"
Set-TextValue $ws.Range("E12") "  +1.74%  "
Set-TextValue $ws.Range("D13") "0.125"
Set-TextValue $ws.Range("E13") "  -2.02%  "
Set-TextValue $ws.Range("D14") "26.64"
Set-TextValue $ws.Range("E14") "  +3.54%  "
Set-TextValue $ws.Range("E15") "  +4.71%  "
Set-TextValue $ws.Range("D16") "56.870.43"
Set-TextValue $ws.Range("E16") "  -0.09%  "
Set-TextValue $ws.Range("B17") "Polkadot"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D17") "6.10"
Set-TextValue $ws.Range("E17") "  +0.41%  "
Set-TextValue $ws.Range("B18") "WrappedEther"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D18") "3.033.52"
Set-TextValue $ws.Range("E18") "  +2.00%  "
Set-TextValue $ws.Range("D19") "13.28"
Set-TextValue $ws.Range("E19") "  +5.40%  "
Set-TextValue $ws.Range("D20") "8.02"
Set-TextValue $ws.Range("E20") "  +3.14%  "
Set-TextValue $ws.Range("D21") "332.76"
Set-TextValue $ws.Range("E21") "  +3.97%  "
Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.11%  "
Set-TextValue $ws.Range("D23") "0.503"
Set-TextValue $ws.Range("E23") "  +3.49%  "
Set-TextValue $ws.Range("D24") "65.32"
Set-TextValue $ws.Range("E24") "  +2.85%  "
Set-TextValue $ws.Range("D25") "3.165.75"
Set-TextValue $ws.Range("E25") "  +1.98%  "
Set-TextValue $ws.Range("E26") "  -0.19%  "
Set-TextValue $ws.Range("E27") "  +1.28%  "
Set-TextValue $ws.Range("D28") "0.0₃0900"
Set-TextValue $ws.Range("E28") "  +1.18%  "
Set-TextValue $ws.Range("D29") "6.66"
Set-TextValue $ws.Range("E29") "  +1.52%  "
Set-TextValue $ws.Range("D30") "7.13"
Set-TextValue $ws.Range("E30") "  +0.21%  "
Set-TextValue $ws.Range("D31") "1.80"
Set-TextValue $ws.Range("E31") "  +2.08%  "
Set-TextValue $ws.Range("E32") "  +3.48%  "
Set-TextValue $ws.Range("D33") "20.60"
Set-TextValue $ws.Range("E33") "  +2.31%  "
Set-TextValue $ws.Range("D34") "4.67"
Set-TextValue $ws.Range("E34") "  +1.03%  "
Set-TextValue $ws.Range("D35") "153.14"
Set-TextValue $ws.Range("E35") "  +0.43%  "
Set-TextValue $ws.Range("D36") "5.91"
Set-TextValue $ws.Range("E36") "  +2.96%  "
Set-TextValue $ws.Range("E37") "  +1.97%  "
Set-TextValue $ws.Range("D38") "25.11"
Set-TextValue $ws.Range("E38") "  +4.68%  "
Set-TextValue $ws.Range("D39") "0.0670"
Set-TextValue $ws.Range("E39") "  +0.68%  "
Set-TextValue $ws.Range("D40") "3.063.06"
Set-TextValue $ws.Range("D41") "36.94"
Set-TextValue $ws.Range("E41") "  -1.55%  "
Set-TextValue $ws.Range("E42") "  +3.62%  "
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  +0.05%  "
Set-TextValue $ws.Range("D44") "0.660"
Set-TextValue $ws.Range("E44") "  +3.31%  "
Set-TextValue $ws.Range("D45") "2.198.27"
Set-TextValue $ws.Range("E45") "  +0.16%  "
Set-TextValue $ws.Range("D46") "1.39"
Set-TextValue $ws.Range("E46") "  +0.88%  "
Set-TextValue $ws.Range("D47") "0.954"
Set-TextValue $ws.Range("E47") "  +0.89%  "
Set-TextValue $ws.Range("E48") "  +0.99%  "
Set-TextValue $ws.Range("D49") "20.18"
Set-TextValue $ws.Range("E49") "  +5.70%  "
Set-TextValue $ws.Range("D50") "0.0242"
Set-TextValue $ws.Range("E50") "  +2.97%  "
Set-TextValue $ws.Range("E51") "  +11.62%  "
